$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, 3, 0, 255, 240, 3, 0, 1, 36),
    @(2, 1, 2, 140, 132, 3, 1, 1, 5),
    @(2, 1, 2, 140, 132, 3, 1, 0, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $data[$i][$j]
    }
}
